$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# Rename the object type labels (replace space with underscore) in column A,
# rows 3-10, to align names of object types.
$ws.Range("A3").Value = "Electrolyzer_AEC"
$ws.Range("A4").Value = "Electrolyzer_PEM"
$ws.Range("A5").Value = "Electrolyzer_SOEC"
$ws.Range("A9").Value = "Methanol_storage"
$ws.Range("A10").Value = "Hydrogen_storage"

# Reflect the selection change recorded for this sheet (A2:A10 selected,
# active cell A2) as part of the editing session.
$ws.Range("A2:A10").Select()
